{"js": "// Add a closing note at the very end of the document body:\n// an empty paragraph, followed by a paragraph containing the note text.\n// Both new paragraphs should pick up the same (\"Normal\"/style0) style\n// that the rest of the document already uses, which is what happens\n// automatically when a new paragraph is inserted after an existing one.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph currently in the document (\"...DVD: /media/haley/HAPPY_FEET\").\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// 1) A blank paragraph right after it.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\n\n// 2) A paragraph with the note text, right after the blank one.\nblankParagraph.insertParagraph(\n  \"* Problems 3 and 5 have screen shots to go along with them if desired\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Add a closing note at the very end of the document body:\n# an empty paragraph, followed by a paragraph containing the note text.\n# New paragraphs inherit the paragraph style (\"Normal\"/style0) already\n# used by the rest of the document, so no explicit style assignment is\n# needed.\n\n$d = $word.ActiveDocument\n\n# The last paragraph currently in the document (\"...DVD: /media/haley/HAPPY_FEET\").\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n# The newly inserted blank paragraph is now the last paragraph.\n$blankParagraph = $d.Paragraphs.Last\n$blankParagraph.Range.InsertParagraphAfter()\n\n# The newly inserted paragraph (now last) gets the note text.\n$noteParagraph = $d.Paragraphs.Last\n$noteParagraph.Range.Text = \"* Problems 3 and 5 have screen shots to go along with them if desired\"\n"}
